$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "total contacts" (L) and "infected x contacts" (M) helper
# columns for each row of the table (m = 0..7 household-contact counts).
$ws.Range("L6").Formula  = "=SUM(B6:K6)"
$ws.Range("M6").Formula  = "=L6*A6"

$ws.Range("L7").Formula  = "=SUM(B7:K7)"
$ws.Range("M7").Formula  = "=L7*A7"

$ws.Range("L8").Formula  = "=SUM(B8:K8)"
$ws.Range("M8").Formula  = "=L8*A8"

$ws.Range("L9").Formula  = "=SUM(B9:K9)"
$ws.Range("M9").Formula  = "=L9*A9"

$ws.Range("L10").Formula = "=SUM(B10:K10)"
$ws.Range("M10").Formula = "=L10*A10"

$ws.Range("L11").Formula = "=SUM(B11:K11)"
$ws.Range("M11").Formula = "=L11*A11"

$ws.Range("L12").Formula = "=SUM(B12:K12)"
$ws.Range("M12").Formula = "=L12*A12"

$ws.Range("L13").Formula = "=SUM(B13:K13)"
$ws.Range("M13").Formula = "=L13*A13"

# Grand total of the M column.
$ws.Range("M14").Formula = "=SUM(M6:M13)"

# Leave the selection where the author ended up.
$ws.Range("O10").Select() | Out-Null
